$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '67.872.24'
Set-TextValue 'E2' '  +1.00%  '

# Row 3
Set-TextValue 'D3' '3.538.80'
Set-TextValue 'E3' '  -0.05%  '

# Row 4
Set-TextValue 'E4' '  +0.05%  '

# Row 5
Set-TextValue 'D5' '615.91'
Set-TextValue 'E5' '  +0.44%  '

# Row 6
Set-TextValue 'D6' '152.47'
Set-TextValue 'E6' '  -1.45%  '

# Row 7
Set-TextValue 'D7' '3.537.93'
Set-TextValue 'E7' '  -0.07%  '

# Row 8
Set-TextValue 'E8' '  -0.03%  '

# Row 9
Set-TextValue 'D9' '0.484'
Set-TextValue 'E9' '  -0.54%  '

# Row 10
Set-TextValue 'D10' '0.140'
Set-TextValue 'E10' '  -1.24%  '

# Row 11
Set-TextValue 'D11' '7.09'
Set-TextValue 'E11' '  +3.19%  '

# Row 12
Set-TextValue 'E12' '  -1.03%  '

# Row 13
Set-TextValue 'D13' '0.0000220'
Set-TextValue 'E13' '  -1.17%  '

# Row 14
Set-TextValue 'D14' '4.139.80'
Set-TextValue 'E14' '  -0.07%  '

# Row 15
Set-TextValue 'D15' '32.14'
Set-TextValue 'E15' '  +0.15%  '

# Row 16
Set-TextValue 'D16' '3.542.17'
Set-TextValue 'E16' '  -0.38%  '

# Row 17
Set-TextValue 'D17' '67.627.32'
Set-TextValue 'E17' '  +0.65%  '

# Row 18
Set-TextValue 'E18' '  -0.58%  '

# Row 19
Set-TextValue 'D19' '6.41'
Set-TextValue 'E19' '  +0.07%  '

# Row 20
Set-TextValue 'D20' '15.38'
Set-TextValue 'E20' '  -0.56%  '

# Row 21
Set-TextValue 'B21' 'Uniswap'
Set-TextValue 'C21' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue 'D21' '9.72'
Set-TextValue 'E21' '  +3.00%  '

# Row 22
Set-TextValue 'B22' 'BitcoinCash'
Set-TextValue 'C22' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 'D22' '447.92'
Set-TextValue 'E22' '  -1.39%  '

# Row 23
Set-TextValue 'E23' '  -2.81%  '

# Row 24
Set-TextValue 'D24' '77.62'
Set-TextValue 'E24' '  -2.30%  '

# Row 25
Set-TextValue 'D25' '0.0000132'
Set-TextValue 'E25' '  +6.06%  '

# Row 26
Set-TextValue 'D26' '3.682.44'
Set-TextValue 'E26' '  -0.01%  '

# Row 27
Set-TextValue 'E27' '  +0.08%  '

# Row 28
Set-TextValue 'D28' '10.27'
Set-TextValue 'E28' '  -1.39%  '

# Row 29
Set-TextValue 'D29' '8.66'
Set-TextValue 'E29' '  +3.19%  '

# Row 30
Set-TextValue 'E30' '  -1.18%  '

# Row 31
Set-TextValue 'D31' '1.60'
Set-TextValue 'E31' '  -4.85%  '

# Row 32
Set-TextValue 'D32' '0.168'
Set-TextValue 'E32' '  +6.82%  '

# Row 33
Set-TextValue 'E33' '  -0.07%  '

# Row 34
Set-TextValue 'D34' '25.98'
Set-TextValue 'E34' '  -0.13%  '

# Row 35
Set-TextValue 'D35' '6.22'
Set-TextValue 'E35' '  +0.27%  '

# Row 36
Set-TextValue 'B36' 'RenzoRestakedETH'
Set-TextValue 'C36' 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue 'D36' '3.528.20'
Set-TextValue 'E36' '  -0.32%  '

# Row 37
Set-TextValue 'B37' 'ImmutableX'
Set-TextValue 'C37' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D37' '1.86'
Set-TextValue 'E37' '  -2.54%  '

# Row 38
Set-TextValue 'D38' '8.06'
Set-TextValue 'E38' '  -0.47%  '

# Row 40
Set-TextValue 'D40' '1.00'
Set-TextValue 'E40' '  +0.03%  '

# Row 41
Set-TextValue 'D41' '176.87'
Set-TextValue 'E41' '  -0.15%  '

# Row 42
Set-TextValue 'D42' '2.20'
Set-TextValue 'E42' '  +2.38%  '

# Row 43
Set-TextValue 'D43' '0.0895'
Set-TextValue 'E43' '  +2.11%  '

# Row 44
Set-TextValue 'D44' '5.43'
Set-TextValue 'E44' '  -3.56%  '

# Row 45
Set-TextValue 'D45' '0.886'
Set-TextValue 'E45' '  -0.92%  '

# Row 46
Set-TextValue 'D46' '28.54'
Set-TextValue 'E46' '  -0.05%  '

# Row 47
Set-TextValue 'D47' '45.38'
Set-TextValue 'E47' '  -1.08%  '

# Row 48
Set-TextValue 'D48' '2.68'
Set-TextValue 'E48' '  -1.11%  '

# Row 49
Set-TextValue 'E49' '  +5.28%  '

# Row 50
Set-TextValue 'D50' '7.63'

# Row 51
Set-TextValue 'D51' '0.996'
Set-TextValue 'E51' '  -4.03%  '
